$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A36").NumberFormat = "@"
$ws.Range("A36").Value = "2025-09-23"
$ws.Range("A36").Style = "Normal"
$ws.Range("B36").Value = "21:20:18"
$ws.Range("C36").Value = "1.00 EUR = 1,629.3841"
